# Apply updated cryptocurrency price (D) and 1h volume-change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "63.705.10"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -3.28%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.607.14"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -2.18%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E4"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "573.55"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -4.43%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "155.99"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -2.45%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E7"; Value = "  +0.07%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "0.624"; ForceText = $true },
    @{ Cell = "E8"; Value = "  -3.13%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "2.604.03"; ForceText = $false },
    @{ Cell = "E9"; Value = "  -2.18%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.118"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -7.20%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "5.82"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -0.99%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "0.381"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -4.85%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.157"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -0.26%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "28.09"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -3.81%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.079.44"; ForceText = $false },
    @{ Cell = "E15"; Value = "  -2.06%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "0.0000180"; ForceText = $true },
    @{ Cell = "E16"; Value = "  -8.00%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "63.570.76"; ForceText = $false },
    @{ Cell = "E17"; Value = "  -3.26%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "2.618.95"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -1.70%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "11.97"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -4.99%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "7.52"; ForceText = $true },
    @{ Cell = "E20"; Value = "  +0.51%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "4.52"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -6.18%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "342.77"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -3.68%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +0.12%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "67.36"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -3.70%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "1.82"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +1.34%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "0.0000108"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -4.62%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "592.96"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +3.82%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "9.12"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -6.41%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "1.57"; ForceText = $true },
    @{ Cell = "E29"; Value = "  -3.49%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "0.161"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -2.01%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E31"; Value = "  +0.09%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "7.90"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -3.10%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "2.05"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -3.60%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "1.73"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -6.36%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "6.56"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -2.63%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "5.41"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -1.69%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.401"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -5.22%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.998"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -0.11%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "19.70"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -4.51%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "154.67"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +0.32%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "1.87"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -5.53%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "0.999"; ForceText = $true },
    @{ Cell = "E42"; Value = "  -0.02%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "2.50"; ForceText = $true },
    @{ Cell = "E43"; Value = "  +1.32%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "41.49"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -3.42%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "156.67"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -3.42%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "23.70"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +0.72%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "3.88"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -5.55%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "0.0589"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -4.71%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "0.628"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -2.60%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "0.1000"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -1.67%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "0.0247"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -4.87%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $r.NumberFormat = "@"
        $r.Value = $u.Value
        $r.Style = "Normal"
    } else {
        $r.Value = $u.Value
    }
}
